$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 31.50994705700387
$ws.Range("C2").Value = 23.72899616659292
$ws.Range("D2").Value = 13.44775367556404
$ws.Range("E2").Value = 13.20295435834256
$ws.Range("G2").Value = 3.828011745734333
$ws.Range("J2").Value = 8.040128651804249
$ws.Range("M2").Value = 28.23008169383322
$ws.Range("N2").Value = 25.97965339208091
$ws.Range("B3").Value = 31.10992175882788
$ws.Range("C3").Value = 23.3713691785962
$ws.Range("D3").Value = 13.45773784624719
$ws.Range("E3").Value = 13.23421245100188
$ws.Range("G3").Value = 3.834600253896437
$ws.Range("J3").Value = 8.05518248244495
$ws.Range("M3").Value = 28.12196706826564
$ws.Range("N3").Value = 25.88680624976296
$ws.Range("B4").Value = 30.87255457891986
$ws.Range("C4").Value = 23.15762935429485
$ws.Range("D4").Value = 13.46669825844024
$ws.Range("E4").Value = 13.25542427850811
$ws.Range("G4").Value = 3.838841804529217
$ws.Range("J4").Value = 8.06488726921196
$ws.Range("M4").Value = 28.06325568213285
$ws.Range("N4").Value = 25.83075370668419
$ws.Range("B5").Value = 30.77800207047485
$ws.Range("C5").Value = 23.07210059053642
$ws.Range("D5").Value = 13.47105903070638
$ws.Range("E5").Value = 13.26457545197659
$ws.Range("G5").Value = 3.840619870551529
$ws.Range("J5").Value = 8.068958737827783
$ws.Range("M5").Value = 28.04126696882299
$ws.Range("N5").Value = 25.80815484898025
$ws.Range("B6").Value = 30.76243603242011
$ws.Range("C6").Value = 23.05799646119329
$ws.Range("D6").Value = 13.47182590367101
$ws.Range("E6").Value = 13.26612561323419
$ws.Range("G6").Value = 3.840918120612536
$ws.Range("J6").Value = 8.069641867312139
$ws.Range("M6").Value = 28.0377329023449
$ws.Range("N6").Value = 25.8044170736242
$ws.Range("B7").Value = 30.87127046707361
$ws.Range("C7").Value = 23.1564693862918
$ws.Range("D7").Value = 13.46675420047902
$ws.Range("E7").Value = 13.2555456416105
$ws.Range("G7").Value = 3.838865582964562
$ws.Range("J7").Value = 8.064941705196764
$ws.Range("M7").Value = 28.06295128496118
$ws.Range("N7").Value = 25.83044794341452
$ws.Range("B8").Value = 31.37036882083151
$ws.Range("C8").Value = 23.60453576197682
$ws.Range("D8").Value = 13.45060749576424
$ws.Range("E8").Value = 13.21331271267665
$ws.Range("G8").Value = 3.830242908673637
$ws.Range("J8").Value = 8.045223753370276
$ws.Range("M8").Value = 28.19121686452676
$ws.Range("N8").Value = 25.94744292538316
$ws.Range("B9").Value = 32.40984736272119
$ws.Range("C9").Value = 24.52492453875453
$ws.Range("D9").Value = 13.44149619504481
$ws.Range("E9").Value = 13.14654229058408
$ws.Range("G9").Value = 3.814878167681477
$ws.Range("J9").Value = 8.010192770757277
$ws.Range("M9").Value = 28.50317237674652
$ws.Range("N9").Value = 26.18433560870756
$ws.Range("B10").Value = 33.20403933700685
$ws.Range("C10").Value = 25.22024154558519
$ws.Range("D10").Value = 13.44867318670993
$ws.Range("E10").Value = 13.10730734266128
$ws.Range("G10").Value = 3.804513836537523
$ws.Range("J10").Value = 7.986633693238923
$ws.Range("M10").Value = 28.76846135192209
$ws.Range("N10").Value = 26.36276714033055
$ws.Range("B11").Value = 33.5705271544378
$ws.Range("C11").Value = 25.53935925933078
$ws.Range("D11").Value = 13.45497193427825
$ws.Range("E11").Value = 13.09159927312968
$ws.Range("G11").Value = 3.799995753165727
$ws.Range("J11").Value = 7.976380742379437
$ws.Range("M11").Value = 28.89678017146815
$ws.Range("N11").Value = 26.44485268882077
$ws.Range("B12").Value = 33.70993485333639
$ws.Range("C12").Value = 25.66049458575582
$ws.Range("D12").Value = 13.45779469507072
$ws.Range("E12").Value = 13.08595947049887
$ws.Range("G12").Value = 3.798312861612448
$ws.Range("J12").Value = 7.97256431277317
$ws.Range("M12").Value = 28.94644854473962
$ws.Range("N12").Value = 26.47606397281861
$ws.Range("B13").Value = 33.67988496546714
$ws.Range("C13").Value = 25.63439468251133
$ws.Range("D13").Value = 13.4571672843799
$ws.Range("E13").Value = 13.0871603700616
$ws.Range("G13").Value = 3.798674061161401
$ws.Range("J13").Value = 7.973383317306064
$ws.Range("M13").Value = 28.93570404218152
$ws.Range("N13").Value = 26.46933647767997
$ws.Range("B14").Value = 33.58198444338056
$ws.Range("C14").Value = 25.5493199098838
$ws.Range("D14").Value = 13.45519538858299
$ws.Range("E14").Value = 13.09112909597197
$ws.Range("G14").Value = 3.799856740926497
$ws.Range("J14").Value = 7.97606544037478
$ws.Range("M14").Value = 28.90084496714841
$ws.Range("N14").Value = 26.44741793576571
$ws.Range("B15").Value = 33.52209556566486
$ws.Range("C15").Value = 25.49724394180422
$ws.Range("D15").Value = 13.45404456311209
$ws.Range("E15").Value = 13.09360025385994
$ws.Range("G15").Value = 3.800584805713727
$ws.Range("J15").Value = 7.97771691346279
$ws.Range("M15").Value = 28.87963234095636
$ws.Range("N15").Value = 26.43400863323738
$ws.Range("B16").Value = 33.18018281030447
$ws.Range("C16").Value = 25.19943342594607
$ws.Range("D16").Value = 13.44832271441548
$ws.Range("E16").Value = 13.1083770310311
$ws.Range("G16").Value = 3.804813038437
$ws.Range("J16").Value = 7.987313036363225
$ws.Range("M16").Value = 28.76022731408327
$ws.Range("N16").Value = 26.35742071500266
$ws.Range("B17").Value = 32.97167676291269
$ws.Range("C17").Value = 25.01737647647222
$ws.Range("D17").Value = 13.44559059501659
$ws.Range("E17").Value = 13.1179907754654
$ws.Range("G17").Value = 3.807457106066675
$ws.Range("J17").Value = 7.993318404749357
$ws.Range("M17").Value = 28.68891773139822
$ws.Range("N17").Value = 26.31066804913311
$ws.Range("B18").Value = 32.85224686920105
$ws.Range("C18").Value = 24.91293364285678
$ws.Range("D18").Value = 13.44430465530729
$ws.Range("E18").Value = 13.12372173834068
$ws.Range("G18").Value = 3.808996435131087
$ws.Range("J18").Value = 7.99681626241151
$ws.Range("M18").Value = 28.64862264224636
$ws.Range("N18").Value = 26.28386366811239
$ws.Range("B19").Value = 32.81189909497036
$ws.Range("C19").Value = 24.87762118986561
$ws.Range("D19").Value = 13.4439182503291
$ws.Range("E19").Value = 13.12569670999661
$ws.Range("G19").Value = 3.809520816672841
$ws.Range("J19").Value = 7.998008105854335
$ws.Range("M19").Value = 28.6351037496594
$ws.Range("N19").Value = 26.27480315668937
$ws.Range("B20").Value = 32.99382200832
$ws.Range("C20").Value = 25.03672947146892
$ws.Range("D20").Value = 13.44585187518994
$ws.Range("E20").Value = 13.11694652724761
$ws.Range("G20").Value = 3.807173724635935
$ws.Range("J20").Value = 7.992674601696543
$ws.Range("M20").Value = 28.69643435931162
$ws.Range("N20").Value = 26.31563602783089
$ws.Range("B21").Value = 33.61072415648489
$ws.Range("C21").Value = 25.57430137997394
$ws.Range("D21").Value = 13.45576269779177
$ws.Range("E21").Value = 13.0899550062421
$ws.Range("G21").Value = 3.799508601284428
$ws.Range("J21").Value = 7.975275845660247
$ws.Range("M21").Value = 28.91105487848782
$ws.Range("N21").Value = 26.45385252560008
$ws.Range("B22").Value = 34.01750284765338
$ws.Range("C22").Value = 25.92728983171152
$ws.Range("D22").Value = 13.46479086613019
$ws.Range("E22").Value = 13.07411302949154
$ws.Range("G22").Value = 3.794662125311367
$ws.Range("J22").Value = 7.964289983787919
$ws.Range("M22").Value = 29.05758703715642
$ws.Range("N22").Value = 26.54492457800132
$ws.Range("B23").Value = 33.80010815397463
$ws.Range("C23").Value = 25.73877828326884
$ws.Range("D23").Value = 13.45973859382278
$ws.Range("E23").Value = 13.08240337238366
$ws.Range("G23").Value = 3.797233946862341
$ws.Range("J23").Value = 7.970118298799719
$ws.Range("M23").Value = 28.97881437952554
$ws.Range("N23").Value = 26.49625157493013
$ws.Range("B24").Value = 32.98380875373452
$ws.Range("C24").Value = 25.0279792757566
$ws.Range("D24").Value = 13.44573286346051
$ws.Range("E24").Value = 13.11741799684279
$ws.Range("G24").Value = 3.807301781495795
$ws.Range("J24").Value = 7.992965523937999
$ws.Range("M24").Value = 28.6930339044659
$ws.Range("N24").Value = 26.31338977274357
$ws.Range("B25").Value = 32.12280369049881
$ws.Range("C25").Value = 24.27211845259561
$ws.Range("D25").Value = 13.44153310095268
$ws.Range("E25").Value = 13.1628834861608
$ws.Range("G25").Value = 3.818871192174816
$ws.Range("J25").Value = 8.019284288512274
$ws.Range("M25").Value = 28.41237834538926
$ws.Range("N25").Value = 26.11946647791462
